# Energy and softmax results
# Fill in the previously empty Energy (E) score row and Softmax aggregation
# values for the BERT / MSP block of the results matrix, then leave the
# active selection on L6 (matching the author's final cursor position).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: BERT / MSP / Softmax
$ws.Range("E7").Value = 76.7
$ws.Range("F7").Value = 49.8
$ws.Range("G7").Value = 8.9
$ws.Range("H7").Value = 27.3

# Row 8: BERT / MSP / E (Energy)
$ws.Range("E8").Value = 78.9
$ws.Range("F8").Value = 54.8
$ws.Range("G8").Value = 17.2
$ws.Range("H8").Value = 26.4

# Move the active selection, matching the saved sheet view state.
$ws.Range("L6").Select()
